$d = $word.ActiveDocument

# Locate the (stable) "_GoBack" bookmark paragraph; the paragraph right
# before it is the empty paragraph that becomes "路由条目的来源", and the
# paragraph right after it is the now-superfluous empty paragraph that gets
# removed once the new content paragraphs are in place.
$bm = $d.Bookmarks.Item("_GoBack")
$bmParaIndex = $bm.Range.Paragraphs.Item(1).Index
$headerParaIndex = $bmParaIndex - 1

# A paragraph that already carries the exact run formatting (theme fonts,
# color, size) we need to reuse for every new run we add.
$templateParaIndex = $headerParaIndex - 1

function Get-TemplateFormattedText() {
    $tp = $d.Paragraphs.Item($templateParaIndex)
    $tr = $d.Range($tp.Range.Start, $tp.Range.End - 1)
    return $tr.FormattedText
}

function Set-ParagraphText($paraIndex, [string]$text) {
    $para = $d.Paragraphs.Item($paraIndex)
    $insertAt = $para.Range.Start
    $collapsed = $d.Range($insertAt, $insertAt)
    $collapsed.FormattedText = (Get-TemplateFormattedText)

    $para2 = $d.Paragraphs.Item($paraIndex)
    $textRange = $d.Range($para2.Range.Start, $para2.Range.End - 1)
    $textRange.Text = $text
}

# 1) Fill the empty paragraph before the bookmark with the section header.
Set-ParagraphText $headerParaIndex "路由条目的来源"

# 2) Insert the six new content paragraphs, each immediately before the
#    (still-empty) bookmark paragraph, filling them in document order.
$lines = @(
    "1.直连路由",
    "路由器本地接口所在的网段",
    "2.静态路由 ",
    "手工配置的路由条目",
    "3.动态路由",
    "路由器之间动态学习到的路由"
)

foreach ($line in $lines) {
    $bmNow = $d.Bookmarks.Item("_GoBack")
    $bmParaNow = $bmNow.Range.Paragraphs.Item(1)
    $bmParaNow.Range.InsertParagraphBefore()

    $bmAfter = $d.Bookmarks.Item("_GoBack")
    $newParaIndex = $bmAfter.Range.Paragraphs.Item(1).Index - 1
    Set-ParagraphText $newParaIndex $line
}

# 3) The paragraph that used to directly follow the bookmark paragraph is
#    now redundant (its content moved into the new paragraphs above) and
#    must be removed, exactly like the diff's deleted empty <w:p>.
$bmFinal = $d.Bookmarks.Item("_GoBack")
$bmParaFinal = $bmFinal.Range.Paragraphs.Item(1)
$trailingIndex = $bmParaFinal.Index + 1
$trailingPara = $d.Paragraphs.Item($trailingIndex)
$trailingPara.Range.Delete()
